{"js": "const body = context.document.body;\n\n// Append a brand-new paragraph at the very end of the document body (right\n// before the final section break), carrying the first sentence as its text.\nconst newParagraph = body.insertParagraph(\n  \"Bot\u00e3o de mudar comparecimento da tabela de consultas n\u00e3o est\u00e1 calculando a diferen\u00e7a de m\u00eas e ano, apenas de dias. \",\n  \"End\"\n);\nnewParagraph.font.italic = false;\nawait context.sync();\n\n// Append the second sentence right after the first one, inside the same\n// paragraph (this lands as its own run immediately following the first).\nconst tailRange = newParagraph.getRange(\"End\");\nconst secondRun = tailRange.insertText(\n  \"Ficar\u00e1 apenas comentado e futuramente ser\u00e1 implementado as configura\u00e7\u00f5es bases do programa, que o usu\u00e1rio poder\u00e1 escolher a quantidade de dias limite para alterar o comparecimento.\",\n  \"End\"\n);\nsecondRun.font.italic = false;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Move to the very end of the document body (just before the final section\n# break) and append a new paragraph there.\n$endRange = $d.Content\n$endRange.Collapse(0)  # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Bot\u00e3o de mudar comparecimento da tabela de consultas n\u00e3o est\u00e1 calculando a diferen\u00e7a de m\u00eas e ano, apenas de dias. \"\n$newPara.Range.Font.Italic = $false\n\n$newPara2 = $d.Paragraphs.Last\n$newPara2.Range.InsertAfter(\"Ficar\u00e1 apenas comentado e futuramente ser\u00e1 implementado as configura\u00e7\u00f5es bases do programa, que o usu\u00e1rio poder\u00e1 escolher a quantidade de dias limite para alterar o comparecimento.\")\n$newPara2.Range.Font.Italic = $false\n"}
